$wb = $excel.ActiveWorkbook

# New values for column B (Transaction Time) for rows 2-4, per sheet
$bValues = @{
    "Todo Results"     = @(54.61, 53.02, 47.46)
    "Project Results"  = @(38.43, 35.42, 28.9)
    "Category Results" = @(31.67, 28.63, 34.05)
}

# Updated values for column C (Create Time) for rows 2-4, per sheet
$cValues = @{
    "Todo Results"     = @(14.78, 14.07, 11.78)
    "Project Results"  = @(13.97, 10.39, 10.11)
    "Category Results" = @(14.47, 10.89, 16.3)
}

foreach ($sheetName in @("Todo Results", "Project Results", "Category Results")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $bVals = $bValues[$sheetName]
    $cVals = $cValues[$sheetName]

    for ($i = 0; $i -lt 3; $i++) {
        $row = $i + 2

        # Add new Transaction Time value in column B
        $ws.Cells.Item($row, 2).Value = $bVals[$i]

        # Update Create Time value in column C
        $ws.Cells.Item($row, 3).Value = $cVals[$i]

        # Remove Delete Time value in column E
        $ws.Cells.Item($row, 5).ClearContents()
    }
}

$wb.Save()
